# remove Gamelogic project, modify SLG building config
#
# The StateFunction sheet is a matrix of which UI "state" (rows, Func1..Func12)
# allows which effect/function (columns B..O). This edit turns off most of the
# effects for each state, keeping only a handful of columns enabled per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..O in row order, 1 = stays enabled, 0 = becomes disabled.
# (Column A is the row's ID / name, untouched.)
$rows = @{
    2  = @(1,1,1,0,0,0,0,0,0,0,0,0,0,1)   # Func1
    3  = @(1,0,0,1,0,0,0,0,0,0,0,0,0,1)   # Func2
    4  = @(1,0,1,0,1,0,0,0,0,0,0,0,0,1)   # Func3
    5  = @(1,0,0,1,0,0,0,0,0,0,0,0,0,1)   # Func4
    6  = @(1,0,0,0,0,0,0,0,0,0,0,0,0,1)   # Func5
    7  = @(1,0,0,0,0,0,0,0,0,0,0,0,0,1)   # Func6
    8  = @(1,0,0,0,0,0,0,0,0,0,0,0,0,1)   # Func7
    9  = @(1,0,0,0,0,0,0,0,0,0,0,0,0,1)   # Func8
    10 = @(1,0,0,0,0,0,0,0,0,0,0,0,0,1)   # Func9
    11 = @(1,0,0,0,0,0,0,0,0,0,0,0,0,1)   # Func10
    12 = @(1,0,0,0,0,0,0,0,0,0,0,0,0,1)   # Func11
    13 = @(1,0,0,0,0,0,0,0,0,0,0,0,0,1)   # Func12
}

foreach ($r in $rows.Keys) {
    $values = $rows[$r]
    for ($i = 0; $i -lt $values.Length; $i++) {
        # column B is index 2
        $col = $i + 2
        $ws.Cells.Item($r, $col).Value = $values[$i]
    }
}

# Move the active selection, as recorded by the save (last place the user clicked).
$ws.Range("F11").Select() | Out-Null
